# Applies the diff: rotates rows 10-13 (row 13's content becomes row 10's,
# row 10's becomes row 11's, row 11's becomes row 12's, row 12's becomes
# row 13's) and appends six brand-new observation rows (14-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to be stored as text even when the string looks like
# a number/date/time (e.g. "1", "2023-09-05", "00:00") so Excel's COM layer
# doesn't silently coerce it into a numeric/date serial value.
function Set-TextCell {
    param($row, $col, $value)
    $cell = $ws.Cells.Item($row, $col)
    if ($value -eq $null -or $value -eq "") {
        $cell.ClearContents()
        return
    }
    $cell.NumberFormat = "@"
    $cell.Value = [string]$value
}

# Helper: plain numeric cell.
function Set-NumCell {
    param($row, $col, $value)
    $ws.Cells.Item($row, $col).Value = $value
}

# Helper: plain boolean cell.
function Set-BoolCell {
    param($row, $col, $value)
    $ws.Cells.Item($row, $col).Value = [bool]$value
}

# Column letter -> index used throughout this sheet.
$COL = @{
    A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8; I=9; J=10; K=11; L=12; M=13; N=14;
    O=15; P=16; Q=17; R=18; S=19; T=20; U=21; V=22; W=23; X=24; Y=25; Z=26;
    AA=27; AB=28; AC=29; AD=30; AE=31; AF=32; AG=33; AH=34; AI=35; AJ=36;
    AK=37; AL=38; AM=39; AN=40; AO=41; AP=42; AQ=43; AR=44; AS=45; AT=46;
    AU=47; AV=48; AW=49; AX=50; AY=51
}

function Write-DataRow {
    param($rowNum, $data)
    foreach ($key in $data.Keys) {
        $col = $COL[$key]
        $val = $data[$key]
        if ($val -is [bool]) {
            Set-BoolCell $rowNum $col $val
        } elseif ($val -is [int] -or $val -is [double]) {
            Set-NumCell $rowNum $col $val
        } else {
            Set-TextCell $rowNum $col $val
        }
    }
}

# ---------------------------------------------------------------------------
# Step 1: rotate the four existing data rows (10-13). Only the cells that
# actually differ between "before" and "after" need touching: A (Id), P
# (Lokalnamn, row 10/11 only), Q (Ost), R (Nord) and AC (Publik kommentar,
# moves from row 11 to row 12).
# ---------------------------------------------------------------------------

# Row 10 <- old row 13's identifying data
Set-NumCell  10 $COL.A 111528980
Set-TextCell 10 $COL.P "Fläcksberget, Hjd"
Set-NumCell  10 $COL.Q 467799.8074815667
Set-NumCell  10 $COL.R 6875539.119922069

# Row 11 <- old row 10's identifying data
Set-NumCell  11 $COL.A 111527876
Set-TextCell 11 $COL.P "Fläcksberget V, Hjd"
Set-NumCell  11 $COL.Q 467615.2905344999
Set-NumCell  11 $COL.R 6875426.740629551
$ws.Cells.Item(11, $COL.AC).ClearContents()

# Row 12 <- old row 11's identifying data
Set-NumCell  12 $COL.A 111528300
Set-NumCell  12 $COL.Q 467795.2212022893
Set-NumCell  12 $COL.R 6875452.272210476
Set-TextCell 12 $COL.AC "Tre blommande."

# Row 13 <- old row 12's identifying data
Set-NumCell  13 $COL.A 111528203
Set-NumCell  13 $COL.Q 467745.6122397452
Set-NumCell  13 $COL.R 6875429.258361855

# ---------------------------------------------------------------------------
# Step 2: append the six new observation rows (14-19).
# ---------------------------------------------------------------------------

Write-DataRow 14 @{
    A=111908768; B=96348; C="Ovaliderad"; D="VU"; E=220787; F="Knärot";
    G="Goodyera repens"; H="(L.) R. Br."; I="1";
    P="Fläcksberget, Hjd"; Q=467911.8445363804; R=6875299.456096188; S=20;
    T="Jämtland"; U="Härjedalen"; V="Härjedalen"; W="Sveg";
    Y="2023-09-05"; Z="00:00"; AA="2023-09-05"; AB="00:00";
    AD=$false; AE=$false; AG=$false;
    AW="lennart karlsson"; AX="lennart karlsson"
}

Write-DataRow 15 @{
    A=111909536; B=77267; C="Ovaliderad"; D="NT"; E=6446; F="Kolflarnlav";
    G="Carbonicola anthracophila"; H="(Nyl.) Bendiksby & Timdal";
    P="Fläcksberget, Hjd"; Q=467891.3929605001; R=6875425.059267788; S=20;
    T="Jämtland"; U="Härjedalen"; V="Härjedalen"; W="Sveg";
    Y="2023-09-05"; Z="00:00"; AA="2023-09-05"; AB="00:00";
    AD=$false; AE=$false; AG=$false;
    AW="lennart karlsson"; AX="lennart karlsson"
}

Write-DataRow 16 @{
    A=111909174; B=77267; C="Ovaliderad"; D="NT"; E=6446; F="Kolflarnlav";
    G="Carbonicola anthracophila"; H="(Nyl.) Bendiksby & Timdal";
    P="Fläcksberget, Hjd"; Q=467989.0228066717; R=6875352.744105402; S=20;
    T="Jämtland"; U="Härjedalen"; V="Härjedalen"; W="Sveg";
    Y="2023-09-05"; Z="00:00"; AA="2023-09-05"; AB="00:00";
    AD=$false; AE=$false; AG=$false;
    AW="lennart karlsson"; AX="lennart karlsson"
}

Write-DataRow 17 @{
    A=111908700; B=89965; C="Ovaliderad"; D="VU"; E=760; F="Doftticka";
    G="Haploporus odorus"; H="(Sommerf.) Bondartsev & Singer"; I="6"; J="fruktkroppar";
    P="Fläcksberget, Hjd"; Q=467921.7931363151; R=6875306.87748003; S=20;
    T="Jämtland"; U="Härjedalen"; V="Härjedalen"; W="Sveg";
    Y="2023-09-05"; Z="00:00"; AA="2023-09-05"; AB="00:00";
    AC="Förekomst av doftticka i avverkningsanmält område.";
    AD=$false; AE=$false; AG=$false;
    AW="lennart karlsson"; AX="lennart karlsson"
}

Write-DataRow 18 @{
    A=111908364; B=90660; C="Ovaliderad"; D="NT"; E=4362; F="Blå taggsvamp";
    G="Hydnellum caeruleum"; H="(Hornem.) P.Karst.";
    P="Gröbäcken, Hjd"; Q=467724.2196293612; R=6874811.291555981; S=20;
    T="Jämtland"; U="Härjedalen"; V="Härjedalen"; W="Sveg";
    Y="2023-09-05"; Z="00:00"; AA="2023-09-05"; AB="00:00";
    AD=$false; AE=$false; AG=$false;
    AW="lennart karlsson"; AX="lennart karlsson"
}

Write-DataRow 19 @{
    A=111909766; B=89183; C="Ovaliderad"; D="LC"; E=3215; F="Rödgul trumpetsvamp";
    G="Craterellus lutescens"; H="(Fr.) Fr.";
    P="Fläcksberget, Hjd"; Q=467756.8135427741; R=6875469.545251801; S=20;
    T="Jämtland"; U="Härjedalen"; V="Härjedalen"; W="Sveg";
    Y="2023-09-05"; Z="00:00"; AA="2023-09-05"; AB="00:00";
    AD=$false; AE=$false; AG=$false;
    AW="lennart karlsson"; AX="lennart karlsson"
}

Write-Host "Done applying edits."
